$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 19.01.2022 08:15"

# Update D5 from text "+0.4" to a real number (0.4)
$ws.Range("D5").Value = 0.4

# Update E5 from text "2022-01-19 08:00:14" to a real date serial value,
# applying the same date/time number format used by the other rows (E2:E10)
$ws.Range("E5").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E5").Value = 44580.33349537037
